$d = $word.ActiveDocument

function Add-EmptyParagraphAfter {
    param($afterIndex)
    $d.Paragraphs($afterIndex).Range.InsertParagraphAfter()
}

function Add-PlainParagraphAfter {
    param($afterIndex, $text)
    $d.Paragraphs($afterIndex).Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs($afterIndex + 1).Range
    $newPara.InsertAfter($text)
}

function Add-LabeledParagraphAfter {
    param($afterIndex, $label, $rest)
    $d.Paragraphs($afterIndex).Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs($afterIndex + 1).Range
    $labelStart = $newPara.Start
    $fullText = $label + $rest
    $newPara.InsertAfter($fullText)
    $labelEnd = $labelStart + $label.Length
    $boldRange = $d.Range($labelStart, $labelEnd)
    $boldRange.Bold = 1
}

# ---------------------------------------------------------------------------
# Step 1: insert two new empty paragraphs before "Análisis Inicial", then
# move "Análisis Inicial" itself below them (so it keeps its own paragraph
# formatting and the two new paragraph marks don't inherit its bold run).
# ---------------------------------------------------------------------------
$d.Paragraphs(2).Range.InsertParagraphBefore()
$d.Paragraphs(2).Range.InsertParagraphBefore()

$d.Paragraphs(1).Range.Cut()

$insertPos = $d.Range($d.Paragraphs(3).Range.Start, $d.Paragraphs(3).Range.Start)
$insertPos.Paste()

# second of the two new leading paragraphs carries a bold paragraph mark
$d.Paragraphs(2).Range.Font.Bold = 1

# ---------------------------------------------------------------------------
# Step 2: cut the "Duración:" paragraph out of its current spot (right after
# "Aplicación seleccionada"); it is re-inserted later, further down.
# ---------------------------------------------------------------------------
$d.Paragraphs(6).Range.Cut()

# ---------------------------------------------------------------------------
# Step 3: after "Se dejó de lado..." (now paragraph 13) add the new
# Equipo / Tester / Metodología lines, then paste "Duración:" back in.
# ---------------------------------------------------------------------------
Add-EmptyParagraphAfter 13
Add-LabeledParagraphAfter 14 "Equipo" ": Ada"
Add-LabeledParagraphAfter 15 "Tester" ": Silvina Vargas"
Add-LabeledParagraphAfter 16 "Metodología" ": Scrum"

$pasteAnchor = $d.Paragraphs(17).Range
$pasteSpot = $d.Range($pasteAnchor.End, $pasteAnchor.End)
$pasteSpot.Paste()

# ---------------------------------------------------------------------------
# Step 4: after the relocated "Duración:" paragraph (now paragraph 18) add
# the Daily Scrum section.
# ---------------------------------------------------------------------------
Add-EmptyParagraphAfter 18
Add-LabeledParagraphAfter 19 "Daily Scrum" ":"
Add-EmptyParagraphAfter 20
Add-PlainParagraphAfter 21 "Día 1: "
Add-PlainParagraphAfter 22 "Tester: Hoy descargaré la app Tiktok lite, me registraré, redactaré historias de usuario e identificaré los criterios de aceptación."
Add-PlainParagraphAfter 23 "Bloqueo: ninguno."
Add-EmptyParagraphAfter 24
Add-PlainParagraphAfter 25 "Día 2:"
Add-PlainParagraphAfter 26 "Tester: haré pruebas funcionales basándome sobre las historias de usuario, identificaré bugs y armaré un reporte de ellos si existieran."
Add-EmptyParagraphAfter 27
Add-PlainParagraphAfter 28 "Día 3: "
Add-PlainParagraphAfter 29 "Tester: Subiré los archivos a un repositorio público para que sean evaluados."
Add-EmptyParagraphAfter 30
Add-PlainParagraphAfter 31 "Día 4: Realizaré pruebas API en Postman, probaré diferentes métodos, haré capturas de pantalla para evidencias. Subiré el archivo a mi repositorio en Github utilizando Git."
Add-EmptyParagraphAfter 32

Write-Output ("final paragraph count=" + $d.Paragraphs.Count)
for ($i = 1; $i -le 36; $i++) {
  Write-Output ($i.ToString() + ": [" + $d.Paragraphs($i).Range.Text + "]")
}
